# Generate Report for Handback
#
# For the "379a23d4-2e39-45f9-81c0-378d661acb84" file row (row 6) on both the
# zh-cn and de-de sheets, the handback generator discovered that the handback
# file it was given is stale (not the latest commit), so it:
#   - records the "Latest Target File" (I6) as a hyperlink back to the
#     source markdown handback doc (same doc A6 already links to),
#   - copies the "Latest Handoff File" name into "Latest Handback File" (J6),
#   - stamps "Latest Handback DateTime" (K6) with the generation time,
#   - and writes an explanatory message into "Error Detail" (P6).
# It also widens the "Error Detail" column (P) so the message is readable.

$wb = $excel.ActiveWorkbook

$handbackMdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b98a0daf1813e78404563521d01dbf1d6eeebaac/e2e/379a23d4-2e39-45f9-81c0-378d661acb84.md"
$handbackMdDisplay = "379a23d4-2e39-45f9-81c0-378d661acb84.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/209de4905d6a7e25cc234dded9f2e9cef641b299/e2e/379a23d4-2e39-45f9-81c0-378d661acb84.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b98a0daf1813e78404563521d01dbf1d6eeebaac/e2e/379a23d4-2e39-45f9-81c0-378d661acb84.md."

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$HandoffFileName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P, the 16th column) so the long
    # message is readable.
    $ws.Range("P:P").ColumnWidth = 39.15

    # I6 "Latest Target File": becomes a hyperlink to the handback markdown
    # file (matching the one A6 already links to).
    $ws.Hyperlinks.Add($ws.Range("I6"), $handbackMdTarget, "", "", $handbackMdDisplay) | Out-Null

    # J6 "Latest Handback File": the xlf file that was actually handed off.
    $ws.Range("J6").Value = $HandoffFileName

    # K6 "Latest Handback DateTime": stamp of when this was generated.
    $ws.Range("K6").Value = $HandbackDateTime

    # P6 "Error Detail": explain why the handback looked stale.
    $ws.Range("P6").Value = $errorDetail
}

Update-LanguageSheet "zh-cn" "379a23d4-2e39-45f9-81c0-378d661acb84.a048c33039853868762fe282aa8ebb863ff7d0bd.zh-cn.xlf" "2016-09-06 20:58:08"

Update-LanguageSheet "de-de" "379a23d4-2e39-45f9-81c0-378d661acb84.a048c33039853868762fe282aa8ebb863ff7d0bd.de-de.xlf" "2016-09-06 20:58:22"
